$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Session 12: add new strategy data for last two sessions
$ws.Range("A13").Value = 12
$ws.Range("C13").Value = 60

# Move the active selection down to the next empty row, matching Excel's
# default behavior after entering data in the row above
$ws.Range("C14").Select()
